$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '44.811.17'
$ws.Range("E2").Value = '  +4.24%  '

# Row 3
$ws.Range("D3").Value = '2.423.70'
$ws.Range("E3").Value = '  +2.71%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.19'
$ws.Range("E5").Value = '  +4.45%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.79'
$ws.Range("E6").Value = '  +6.88%  '

# Row 7
$ws.Range("E7").Value = '  +2.63%  '

# Row 8
$ws.Range("E8").Value = '  -0.05%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.532'
$ws.Range("E9").Value = '  +11.64%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.43'
$ws.Range("E10").Value = '  +3.13%  '

# Row 11
$ws.Range("E11").Value = '  +1.82%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.70'
$ws.Range("E12").Value = '  +1.22%  '

# Row 13
$ws.Range("E13").Value = '  -1.74%  '

# Row 14
$ws.Range("E14").Value = '  +3.50%  '

# Row 15
$ws.Range("D15").Value = '2.802.11'
$ws.Range("E15").Value = '  +2.79%  '

# Row 16
$ws.Range("D16").Value = '2.389.79'
$ws.Range("E16").Value = '  +1.12%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.835'
$ws.Range("E17").Value = '  +4.90%  '

# Row 18
$ws.Range("D18").Value = '44.639.80'
$ws.Range("E18").Value = '  +3.86%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.33'
$ws.Range("E19").Value = '  +3.18%  '

# Row 20
$ws.Range("E20").Value = '  +1.93%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0917'
$ws.Range("E21").Value = '  +3.56%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.75'
$ws.Range("E22").Value = '  +1.21%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '242.46'
$ws.Range("E23").Value = '  +3.11%  '

# Row 24
$ws.Range("E24").Value = '  +4.22%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.49'
$ws.Range("E25").Value = '  +2.04%  '

# Row 26
$ws.Range("E26").Value = '  -0.08%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.19'
$ws.Range("E27").Value = '  +3.34%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.28'
$ws.Range("E28").Value = '  -3.79%  '

# Row 29
$ws.Range("E29").Value = '  +1.50%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.70'
$ws.Range("E30").Value = '  +4.17%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '48.64'
$ws.Range("E31").Value = '  +1.75%  '

# Row 32
$ws.Range("E32").Value = '  +18.82%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.57'
$ws.Range("E33").Value = '  +11.75%  '

# Row 34
$ws.Range("E34").Value = '  +3.16%  '

# Row 35
$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.25%  '

# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0769'
$ws.Range("E36").Value = '  +6.07%  '

# Row 37
$ws.Range("E37").Value = '  +3.75%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.48'
$ws.Range("E38").Value = '  +3.59%  '

# Row 39
$ws.Range("E39").Value = '  +0.79%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '123.53'
$ws.Range("E40").Value = '  -3.69%  '

# Row 41
$ws.Range("E41").Value = '  +1.87%  '

# Row 42
$ws.Range("E42").Value = '  -2.96%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.94'
$ws.Range("E43").Value = '  +0.22%  '

# Row 44
$ws.Range("E44").Value = '  +4.12%  '

# Row 45
$ws.Range("D45").Value = '1.943.41'
$ws.Range("E45").Value = '  +0.85%  '

# Row 46
$ws.Range("E46").Value = '  +8.09%  '

# Row 47
$ws.Range("E47").Value = '  -2.39%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.30'
$ws.Range("E48").Value = '  +0.53%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.76'
$ws.Range("E49").Value = '  +16.49%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '75.94'
$ws.Range("E50").Value = '  +6.37%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '53.99'
$ws.Range("E51").Value = '  +5.42%  '
